# Fix Training Data Issue
# The BF column holds a "Date" label (row 1) followed by a date string for each
# team row. The stored text "5-1-2012-13" was off by one day (NBA stats quirk),
# so it needs to become the correct ISO-style date string "2013-05-01".
#
# Column BF = column index 58.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "5-1-2012-13"
$newText = "2013-05-01"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # Column BF
    if ($cell.Value2 -eq $oldText) {
        # Assigning the literal date-like string directly via Value2/Value would
        # make Excel re-interpret it as a real date (serial number). To keep it
        # as plain text matching the target, write it through a text formula and
        # then convert that formula result to a static value via copy/paste.
        $cell.Formula = '="' + $newText + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

$excel.CutCopyMode = 0
